$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-run simulation (RU 1001, without crop): a new "Empty" category is
# inserted at the top of the keyword table and every existing category
# shifts down one row; "Corruption" (previously the last populated row)
# moves into the row that used to be the all-zero "Empty" placeholder.

$labels = @("Empty", "Nothing", "Compliment", "Interesting", "Thank", "Difficult", "Corruption")
for ($i = 0; $i -lt $labels.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $labels[$i]
}

$data = @(
    @(0.474382958912649, 0.661465379406704, 0.619213855609162, 0.718653077904179, 0.699123724855438, 0.568263429037403, 0.645246378200891, 0.657384872305653, 0.699072247882384, 0.414192850264813, 0, 0.402042003578956, 0.493972753411824),
    @(0.076694105499444, 0.0323173199926413, 0.0271393975234657, 0.0119858594027523, 0.0209287217995857, 0.0797427421450919, 0.0328774847829057, 0.048330321078416, 0.0143718808059912, 0.220305312570907, 0.0460794745575303, 0.0548672526707154, 0.0854953417813122),
    @(0.0584846703368328, 0.0385153922965894, 0.0448715527200063, 0.0285939348605774, 0.0249686958339627, 0.0261593011672415, 0.0408290728299135, 0.0630478921225709, 0.0366532821438346, 0.0239311513840183, 0.0525084921753547, 0.13117776919282, 0.0909565818320698),
    @(0.0531584284123687, 0.0509002944878847, 0.0602364423909903, 0.0279149677024095, 0.0596263552340503, 0.0597336422893675, 0.0614955273531362, 0.0506657617476775, 0.0453220880472148, 0.0280269954650683, 0.0407140857166648, 0.0185888292630765, 0.0739803017703995),
    @(0.0338923717916661, 0.025331809777643, 0.0257473393482286, 0.0267882954995236, 0.0175195795419402, 0.0351048342123718, 0.0316774452370636, 0.0208642635125292, 0.0179095670244576, 0.0048953074515604, 0.0277091539062055, 0.13122521787311, 0.0491228741835903),
    @(0.0197917773969517, 0.0109296860457397, 0.0142817379883159, 0.00580153130414161, 0.00925101670355077, 0.00728772658819937, 0.0306761015195065, 0.00391829685666812, 0.00588811910333217, 0.0842519919712641, 0.0140821167599252, 0.0013074277785911, 0.00851440446614129),
    @(0.00392784933537338, 0.00398359150418561, 0.00367194117586363, 0.00506433468696284, 0.00247983796047947, 0.00189755228820627, 0.0056971952956279, 0.00329082456419541, 0.0111732486072152, 0, 0.00615643692764775, 0.00333590233748943, 0.00451531276518934)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 2).Value = $rowValues[$c]
    }
}
